$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text (string) type instead of
# being auto-converted to numbers/dates/percentages by Excel.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '53.604.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  -5.10%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.239.54'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -5.90%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.998'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.28%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '486.94'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -3.83%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '125.35'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -4.26%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.21%  '

$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -4.34%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.237.06'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -6.45%  '

$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -6.88%  '

$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -0.26%  '

$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -2.89%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.63'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -4.55%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.607.94'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -6.99%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.24'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -1.94%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '53.529.75'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -5.18%  '

$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -4.46%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.222.05'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  -7.35%  '

$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -4.47%  '

$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.22%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '295.61'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -4.71%  '

$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.45%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.995'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.43%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.66'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -4.25%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.996'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -0.02%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.368'

$ws.Range("B27").NumberFormat = "@"
$ws.Range("B27").Value = 'Kaspa'
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.147'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -0.62%  '

$ws.Range("B28").NumberFormat = "@"
$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("C28").NumberFormat = "@"
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.315.70'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -7.15%  '

$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -3.18%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '163.82'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -5.49%  '

$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -4.20%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.997'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.16%  '

$ws.Range("B33").NumberFormat = "@"
$ws.Range("B33").Value = 'PEPE'
$ws.Range("C33").NumberFormat = "@"
$ws.Range("C33").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0₃0669'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -6.51%  '

$ws.Range("B34").NumberFormat = "@"
$ws.Range("B34").Value = 'Aptos'
$ws.Range("C34").NumberFormat = "@"
$ws.Range("C34").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.78'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -1.32%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.996'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.02%  '

$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.42%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '17.28'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -2.52%  '

$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.25%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.832'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  -4.11%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '35.17'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -3.72%  '

$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -1.06%  '

$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.46%  '

$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -2.99%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '126.54'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -1.68%  '

$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.97%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0882'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.91%  '

$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -5.56%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '234.64'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -2.57%  '

$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -2.26%  '

$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.27%  '
